$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case row for the "myMeth_test" Test table (row 38 was the last one)
$ws.Range("B39").Value = 1
$ws.Range("C39").Value = 1

# Update existing myDataSr row: subField.field2 now shows combined array values
$ws.Range("D45").Value = "test33, test44"

# New data row appended to the myDataSr data table
$ws.Range("B46").Value = 1
$ws.Range("C46").Value = "text1"
$ws.Range("D46").Value = "5, 2"
$ws.Range("E46").Value = "aa"
$ws.Range("F46").Value = "aaa"
$ws.Range("G46").Value = "bb"
$ws.Range("H46").Value = "bbb"

# Update the selection / viewport to match the edited area
$ws.Range("C39").Select()
$ws.Application.ActiveWindow.ScrollRow = 10
